## Update US model copy
## - About sheet: refresh source figures (value + inflation factor), label year,
##   turn the source-link cell into a real hyperlink, tidy selection/row height.
## - CCSTaSC sheet: recalculates off the About inputs; clear the stray explicit
##   number-format override on the result row so it inherits the default style.

$wb  = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$ccs   = $wb.Worksheets.Item("CCSTaSC")

# --- About sheet updates ---------------------------------------------------

# Updated source figures.
$about.Range("B10").Value = 23
$about.Range("B11").Value = 0.9143273584567535

# Relabel the inflation-adjustment row.
$about.Range("A11").Value = "2018 to 2012 $"

# Turn the source URL text into a clickable hyperlink (keeps existing text
# and formatting - Hyperlinks.Add likes to stamp an extra explicit font
# attribute on the cell, so stash the original format and restore it after).
$srcCell = $about.Range("B7")
$srcUrl = "https://netzeroamerica.princeton.edu/img/Princeton%20NZA%20FINAL%20REPORT%20SUMMARY%20(29Oct2021).pdf"
$scratch = $about.Range("D1")
$srcCell.Copy()
$scratch.PasteSpecial(-4122)  # xlPasteFormats
$srcCell.Hyperlinks.Add($srcCell, $srcUrl)
$scratch.Copy()
$srcCell.PasteSpecial(-4122)  # xlPasteFormats
$scratch.Clear()

# The wrapped source-link text now renders across more lines, so the row
# grows taller to fit it.
$about.Rows.Item(7).RowHeight = 45

# --- CCSTaSC sheet updates --------------------------------------------------

# The result row no longer needs its explicit number-format override -
# clearing formats lets it fall back to the default style (and drops the
# now-unused cell format from the workbook).
$ccs.Range("B2:AE2").ClearFormats()

# --- Selection / activation bookkeeping ------------------------------------

[void]$ccs.Range("B2").Select()
[void]$ccs.Activate()
[void]$about.Range("B11").Select()
[void]$about.Activate()

$wb.Application.Calculate()
